$p = $ppt.ActivePresentation

foreach ($idx in 14,15,16) {
    $s = $p.Slides.Item($idx)
    $shp = $s.Shapes.Item(1)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{0EE82FCD-CF5D-45DF-9DC5-F85D7F1B6AF9}")
    }
}
